$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.169212666666667
$ws.Range("H2").Value = 3.507638
$ws.Range("I2").Value = 0.0005193657195729173
$ws.Range("J2").Value = 0.0005193657195729173
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.121657333333333
$ws.Range("N2").Value = 3.364972
$ws.Range("O2").Value = 0.01078859740606296
$ws.Range("P2").Value = 0.01078859740606296
$ws.Range("Q2").Value = 1.311455961792889
$ws.Range("R2").Value = 11.803103656136
$ws.Range("S2").Value = 0.0000056032276549823969991059868789307075
$ws.Range("T2").Value = 0.0000056032276549823969991059868789307075
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.169212666666667
$ws.Range("H3").Value = 3.507638
$ws.Range("I3").Value = 0.0005193657195729173
$ws.Range("J3").Value = 0.0005193657195729173
$ws.Range("O3").Value = 0.8063022765396375
$ws.Range("P3").Value = 0.8063022765396375
$ws.Range("Q3").Value = 98.01366088429934
$ws.Range("R3").Value = 882.122947958694
$ws.Range("S3").Value = 0.0004187657620482902
$ws.Range("T3").Value = 0.0004187657620482902
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.169212666666667
$ws.Range("H4").Value = 3.507638
$ws.Range("I4").Value = 0.0005193657195729173
$ws.Range("J4").Value = 0.0005193657195729173
$ws.Range("M4").Value = 19.0165
$ws.Range("N4").Value = 57.0495
$ws.Range("O4").Value = 0.1829091260542996
$ws.Range("P4").Value = 0.1829091260542996
$ws.Range("Q4").Value = 22.23433267566667
$ws.Range("R4").Value = 200.108994081
$ws.Range("S4").Value = 0.0000949967298696447669205042418916207225
$ws.Range("T4").Value = 0.0000949967298696447669205042418916207225
$ws.Range("I5").Value = 0.9638330474556795
$ws.Range("J5").Value = 0.9638330474556795
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.121657333333333
$ws.Range("N5").Value = 3.364972
$ws.Range("O5").Value = 0.01078859740606296
$ws.Range("P5").Value = 0.01078859740606296
$ws.Range("Q5").Value = 2433.785189554264
$ws.Range("R5").Value = 21904.06670598837
$ws.Range("S5").Value = 0.0103984067156581
$ws.Range("T5").Value = 0.0103984067156581
$ws.Range("I6").Value = 0.9638330474556795
$ws.Range("J6").Value = 0.9638330474556795
$ws.Range("O6").Value = 0.8063022765396375
$ws.Range("P6").Value = 0.8063022765396375
$ws.Range("S6").Value = 0.7771407803676509
$ws.Range("T6").Value = 0.7771407803676509
$ws.Range("I7").Value = 0.9638330474556795
$ws.Range("J7").Value = 0.9638330474556795
$ws.Range("M7").Value = 19.0165
$ws.Range("N7").Value = 57.0495
$ws.Range("O7").Value = 0.1829091260542996
$ws.Range("P7").Value = 0.1829091260542996
$ws.Range("Q7").Value = 41262.223926819
$ws.Range("R7").Value = 371360.015341371
$ws.Range("S7").Value = 0.1762938603723706
$ws.Range("T7").Value = 0.1762938603723706
$ws.Range("G8").Value = 80.250984
$ws.Range("H8").Value = 240.752952
$ws.Range("I8").Value = 0.03564758682474761
$ws.Range("J8").Value = 0.0356475868247476
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.121657333333333
$ws.Range("N8").Value = 3.364972
$ws.Range("O8").Value = 0.01078859740606296
$ws.Range("P8").Value = 0.01078859740606296
$ws.Range("Q8").Value = 90.014104710816
$ws.Range("R8").Value = 810.1269423973439
$ws.Range("S8").Value = 0.0003845874627498761
$ws.Range("T8").Value = 0.0003845874627498759
$ws.Range("G9").Value = 80.250984
$ws.Range("H9").Value = 240.752952
$ws.Range("I9").Value = 0.03564758682474761
$ws.Range("J9").Value = 0.0356475868247476
$ws.Range("O9").Value = 0.8063022765396375
$ws.Range("P9").Value = 0.8063022765396375
$ws.Range("Q9").Value = 6727.341360260664
$ws.Range("R9").Value = 60546.07224234597
$ws.Range("S9").Value = 0.02874273040993838
$ws.Range("T9").Value = 0.02874273040993838
$ws.Range("G10").Value = 80.250984
$ws.Range("H10").Value = 240.752952
$ws.Range("I10").Value = 0.03564758682474761
$ws.Range("J10").Value = 0.0356475868247476
$ws.Range("M10").Value = 19.0165
$ws.Range("N10").Value = 57.0495
$ws.Range("O10").Value = 0.1829091260542996
$ws.Range("P10").Value = 0.1829091260542996
$ws.Range("Q10").Value = 1526.092837236
$ws.Range("R10").Value = 13734.835535124
$ws.Range("S10").Value = 0.006520268952059351
$ws.Range("T10").Value = 0.006520268952059349
